$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Character")

# Add a new character row (row 6), copying formatting from the row above (row 5)
# so the new row picks up the same cell style (border) instead of creating a
# brand-new style entry.
$ws.Range("A5:H5").Copy($ws.Range("A6:H6"))

$ws.Range("A6").Value = 2
$ws.Range("B6").Value = "chr_bird"
$ws.Range("C6").Value = 400
$ws.Range("D6").Value = 100
$ws.Range("E6").Value = 200
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2
$ws.Range("H6").Value = 3

# Make the Character sheet the active tab and select the newly added cell,
# matching the popup-system edit moving focus away from the Gimmick sheet.
$ws.Activate()
$ws.Range("E6").Select()
